$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the contact e-mail address (cell B8) -------------------------
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com"

# Re-point the existing mailto: hyperlink on B8 at the new address, and
# drop its stale cached display text so Excel regenerates it from the cell.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$8') {
        $h.Address = "mailto:yryskan.kalymbetova@gmail.com"
        $h.TextToDisplay = ""
    }
}

# --- Add a hyperlink for the organisation web site (cell B10) ------------
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.stat.gov.kg/")

# Match the existing "hyperlink" cell formatting used by B8 (same visual
# style as the other contact-details hyperlink in this sheet).
$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Move the active selection from B19 to B9 -----------------------------
$ws.Range("B9").Select()
